# Generate Report for Handoff
# Updates the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps for the
# fa66597a-1734-4442-bcfd-ea0286d426a7.md file after a new handoff was generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# fa66597a-1734-4442-bcfd-ea0286d426a7.md row (row 5)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-08-31 07:17:42"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the
# fa66597a-1734-4442-bcfd-ea0286d426a7.md row (row 5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-08-31 07:17:31"

# de-de sheet: "Latest Handoff Datetime" column (H) for the
# fa66597a-1734-4442-bcfd-ea0286d426a7.md row (row 5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-08-31 07:17:42"
